# Update values in column C ("Name of Algo" result column) for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C11" = -13.08749999999999
    "C12" = -11.3877
    "C15" = -13.47869999999999
    "C27" = -12.968
    "C28" = -13.0874
    "C31" = -13.43290000000001
    "C32" = -13.85750000000001
    "C36" = -12.50390000000001
    "C38" = -12.40009999999999
    "C46" = -14.5423
    "C54" = -13.05070000000001
    "C55" = -14.03600000000001
    "C56" = -12.4944
    "C67" = -10.82110000000001
    "C69" = -12.13559999999999
    "C72" = -11.3516
    "C73" = -12.30300000000001
    "C83" = -13.7116
    "C86" = -13.49449999999999
    "C91" = -10.3075
    "C93" = -10.7436
    "C99" = -12.87749999999999
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
